$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @(
  '74-65=',
  '6+38=',
  '19+53=',
  '40-28=',
  '92-89=',
  '46-18=',
  '48+38=',
  '52-29=',
  '32-6=',
  '47+5=',
  '6+36=',
  '45-19=',
  '36-29=',
  '65-28=',
  '66+6=',
  '65-16=',
  '19+58=',
  '62-14=',
  '5+7=',
  '74-26=',
  '83-49=',
  '47+34=',
  '62-19=',
  '49+5=',
  '22+49=',
  '65-49=',
  '42-39=',
  '25+37=',
  '94-47=',
  '8+7=',
  '8+37=',
  '80-32=',
  '91-22=',
  '50-32=',
  '46-8=',
  '49+14=',
  '44-18=',
  '72-58=',
  '66-19=',
  '25+38=',
  '37+27=',
  '12-6=',
  '93-24=',
  '84-18=',
  '3+48=',
  '27+18=',
  '29+43=',
  '91-58=',
  '18+37=',
  '97-79=',
  '90-43=',
  '8+18=',
  '39+12=',
  '29+19=',
  '31-15=',
  '60-4=',
  '39+13=',
  '89+7=',
  '70-31=',
  '72-26=',
  '70-47=',
  '70-64=',
  '54-35=',
  '83-19=',
  '19+2=',
  '82-7=',
  '80-73=',
  '79+9=',
  '66+27=',
  '94-85=',
  '87-58=',
  '93-36=',
  '5+29=',
  '91-3=',
  '49+43=',
  '68+27=',
  '58+37=',
  '57+26=',
  '29+66=',
  '64-48=',
  '44-39=',
  '51-8=',
  '51-24=',
  '24-6=',
  '73-26=',
  '18+28=',
  '32-4=',
  '79+17=',
  '17+16=',
  '88+4=',
  '76-37=',
  '63-9=',
  '86-37=',
  '3+59=',
  '41-17=',
  '36+19=',
  '49+34=',
  '71-14=',
  '37+17=',
  '8+57='
)
$cols = 5
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = [int][math]::Floor($i / $cols) + 1
    $col = ($i % $cols) + 1
    $t.Cell($row, $col).Range.Text = $values[$i]
}
Write-Output "done"
